$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns D (Price) and E (Volume) remain text, not auto-converted numbers
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.231.95"
$ws.Range("E2").Value = "  -0.76%  "

$ws.Range("D3").Value = "1.874.44"
$ws.Range("E3").Value = "  -1.93%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "234.71"
$ws.Range("E5").Value = "  -1.83%  "

$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Value = "0.4689"
$ws.Range("E7").Value = "  -1.95%  "

$ws.Range("D8").Value = "0.2835"
$ws.Range("E8").Value = "  -0.32%  "

$ws.Range("D9").Value = "0.06593"
$ws.Range("E9").Value = "  -1.78%  "

$ws.Range("D10").Value = "20.32"
$ws.Range("E10").Value = "  +7.88%  "

$ws.Range("D11").Value = "0.07777"
$ws.Range("E11").Value = "  +0.98%  "

$ws.Range("D12").Value = "97.59"
$ws.Range("E12").Value = "  -4.54%  "

$ws.Range("D13").Value = "1.883.50"
$ws.Range("E13").Value = "  -1.53%  "

$ws.Range("D14").Value = "5.078"
$ws.Range("E14").Value = "  -2.58%  "

$ws.Range("D15").Value = "0.6721"
$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("D16").Value = "284.75"
$ws.Range("E16").Value = "  +6.48%  "

$ws.Range("D17").Value = "30.262.28"
$ws.Range("E17").Value = "  -0.76%  "

$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").Value = "12.61"
$ws.Range("E19").Value = "  -0.61%  "

$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.121.57"
$ws.Range("E20").Value = "  -1.72%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "5.386"
$ws.Range("E21").Value = "  -0.51%  "

$ws.Range("D22").Value = "0.000007264"
$ws.Range("E22").Value = "  -2.83%  "

$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "6.165"
$ws.Range("E24").Value = "  -2.10%  "

$ws.Range("D25").Value = "9.369"
$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").Value = "167.20"
$ws.Range("E26").Value = "  +0.24%  "

$ws.Range("D27").Value = "19.18"
$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("D28").Value = "1.977"
$ws.Range("E28").Value = "  -4.29%  "

$ws.Range("D29").Value = "1.371"
$ws.Range("E29").Value = "  -0.86%  "

$ws.Range("D30").Value = "0.09668"
$ws.Range("E30").Value = "  -3.47%  "

$ws.Range("D31").Value = "4.382"
$ws.Range("E31").Value = "  -5.33%  "

$ws.Range("D32").Value = "1.469"
$ws.Range("E32").Value = "  -2.83%  "

$ws.Range("D33").Value = "4.112"
$ws.Range("E33").Value = "  -2.70%  "

$ws.Range("D34").Value = "0.04686"
$ws.Range("E34").Value = "  -1.02%  "

$ws.Range("D35").Value = "0.7061"
$ws.Range("E35").Value = "  -2.81%  "

$ws.Range("D36").Value = "1.092"
$ws.Range("E36").Value = "  -1.43%  "

$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "0.9999"
$ws.Range("E37").Value = "  +0.12%  "

$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "2.717"
$ws.Range("E38").Value = "  +0.06%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01869"
$ws.Range("E39").Value = "  -2.44%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "6.481"
$ws.Range("E40").Value = "  +3.14%  "

$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "2.524"
$ws.Range("E41").Value = "  -3.31%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "71.96"
$ws.Range("E42").Value = "  -4.43%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "1.955"
$ws.Range("E43").Value = "  -0.76%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "0.8609"
$ws.Range("E44").Value = "  +0.11%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "102.99"
$ws.Range("E46").Value = "  -1.87%  "

$ws.Range("B47").Value = "TheSandbox"
$ws.Range("C47").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D47").Value = "0.4187"
$ws.Range("E47").Value = "  -1.83%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "984.69"
$ws.Range("E48").Value = "  +7.08%  "

$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "7.215"
$ws.Range("E49").Value = "  -2.83%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.125"
$ws.Range("E50").Value = "  +3.83%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "33.94"
$ws.Range("E51").Value = "  -2.54%  "
